$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range('D2').Value = '69.582.65'
$ws.Range('E2').Value = '  -1.65%  '

$ws.Range('D3').Value = '3.497.86'
$ws.Range('E3').Value = '  -1.69%  '

$ws.Range('E4').Value = '  -0.16%  '

$ws.Range('D5').Value = '611.36'
$ws.Range('E5').Value = '  +5.19%  '

$ws.Range('D6').Value = '188.19'
$ws.Range('E6').Value = '  +0.56%  '

$ws.Range('D7').Value = '0.625'
$ws.Range('E7').Value = '  -0.28%  '

$ws.Range('D8').Value = '0.998'
$ws.Range('E8').Value = '  -0.25%  '

$ws.Range('D9').Value = '0.214'
$ws.Range('E9').Value = '  -3.80%  '

$ws.Range('D10').Value = '0.648'
$ws.Range('E10').Value = '  -0.08%  '

$ws.Range('D11').Value = '53.09'
$ws.Range('E11').Value = '  -2.79%  '

$ws.Range('E12').Value = '  -3.74%  '

$ws.Range('E13').Value = '  +0.29%  '

$ws.Range('D14').Value = '4.061.70'
$ws.Range('E14').Value = '  -1.63%  '

$ws.Range('D15').Value = '599.09'
$ws.Range('E15').Value = '  +4.95%  '

$ws.Range('D16').Value = '69.656.15'
$ws.Range('E16').Value = '  -1.65%  '

$ws.Range('D17').Value = '18.95'
$ws.Range('E17').Value = '  -1.11%  '

$ws.Range('D18').Value = '12.60'
$ws.Range('E18').Value = '  -1.80%  '

$ws.Range('D19').Value = '3.492.53'
$ws.Range('E19').Value = '  -2.49%  '

$ws.Range('E20').Value = '  -0.30%  '

$ws.Range('E21').Value = '  -1.32%  '

$ws.Range('D22').Value = '17.28'
$ws.Range('E22').Value = '  -2.41%  '

$ws.Range('D23').Value = '105.19'
$ws.Range('E23').Value = '  +12.33%  '

$ws.Range('D24').Value = '5.11'
$ws.Range('E24').Value = '  +4.57%  '

$ws.Range('D25').Value = '4.66'
$ws.Range('E25').Value = '  +2.51%  '

$ws.Range('D26').Value = '3.06'
$ws.Range('E26').Value = '  +3.55%  '

$ws.Range('E27').Value = '  -2.83%  '

$ws.Range('D28').Value = '9.75'
$ws.Range('E28').Value = '  +5.59%  '

$ws.Range('D29').Value = '33.42'
$ws.Range('E29').Value = '  +3.15%  '

$ws.Range('D30').Value = '6.96'
$ws.Range('E30').Value = '  -3.36%  '

$ws.Range('D31').Value = '4.13'
$ws.Range('E31').Value = '  +13.66%  '

$ws.Range('D32').Value = '12.46'
$ws.Range('E32').Value = '  +1.40%  '

$ws.Range('E33').Value = '  -0.85%  '

$ws.Range('D34').Value = '63.51'
$ws.Range('E34').Value = '  +0.71%  '

$ws.Range('D35').Value = '3.16'
$ws.Range('E35').Value = '  -6.04%  '

$ws.Range('D36').Value = '0.999'
$ws.Range('E36').Value = '  -0.05%  '

$ws.Range('E37').Value = '  +8.21%  '

$ws.Range('D38').Value = '515.14'
$ws.Range('E38').Value = '  -4.77%  '

$ws.Range('B39').Value = 'Maker'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D39').Value = '3.613.76'
$ws.Range('E39').Value = '  +1.16%  '

$ws.Range('B40').Value = 'TheGraph'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D40').Value = '0.395'
$ws.Range('E40').Value = '  -4.30%  '

$ws.Range('D41').Value = '36.81'
$ws.Range('E41').Value = '  -3.65%  '

$ws.Range('D42').Value = '0.0₃0776'
$ws.Range('E42').Value = '  -3.11%  '

$ws.Range('E43').Value = '  -0.81%  '

$ws.Range('E44').Value = '  -0.76%  '

$ws.Range('E45').Value = '  +0.01%  '

$ws.Range('E46').Value = '  +2.84%  '

$ws.Range('D47').Value = '3.36'
$ws.Range('E47').Value = '  -3.32%  '

$ws.Range('D48').Value = '8.79'
$ws.Range('E48').Value = '  -5.68%  '

$ws.Range('E49').Value = '  +0.30%  '

$ws.Range('D50').Value = '131.54'
$ws.Range('E50').Value = '  -1.85%  '

$ws.Range('E51').Value = '  -7.93%  '
